$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The importer used to write a placeholder/duplicate row (row 2) before the
# real data rows. Remove it so the real data starts on row 2.
$ws.Rows("2:2").Delete()

# Re-apply the (non-bold) default style so the moved "MATERIA_PRIMA" cell no
# longer drags along the stale "applyFont" style it had before.
$ws.Range("C3").Font.Bold = $false

# Make sure the date column keeps its mm/yyyy display format.
$ws.Range("E2:E6").NumberFormat = "mm/yyyy"

# Reflect where the user left the selection after fixing the sheet
# (an entire-row selection on the new first data row).
$null = $ws.Range("A2:XFD2").Select()
